$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GPLIM-2588 Fix spreadsheet headers.
# A1 keeps its bold/quote-prefixed header style, so re-enter the text with a
# leading apostrophe (forces text entry without altering the cell's xf/style,
# mirroring how the cell was originally authored).
$ws.Range("A1").Value = "'Specimen_Number"
$ws.Range("F1").Value = "SAMPLE_TYPE"

# Leave the selection on the last-edited header cell.
$ws.Range("F1").Select()
